$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text cells keep their literal string representation (no numeric/date coercion).

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.862.37"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.20%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.105.19"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.28%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "525.28"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.30%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.35"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.21%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.08%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.104.94"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.31%  "

# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.65%  "

# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.13%  "

# Row 11
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.36%  "

# Row 12
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +3.96%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.638.35"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.24%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.133"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.98%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.64"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -3.24%  "

# Row 16
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.25%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "57.959.78"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.25%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.108.49"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.42%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.09"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.03%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.77"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.16%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.99"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.76%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "341.46"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.88%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.17%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.513"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.56%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "67.23"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +3.74%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.170"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.24%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.12%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0918"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.84%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.46"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.56%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.27"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.59%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.87"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +4.32%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.01"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.19%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.20"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.90%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "158.07"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.81%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.66"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +2.33%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.17"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.91%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "26.42"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -4.27%  "

# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.60%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0669"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.75%  "

# Row 41
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.58"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +14.38%  "

# Row 42
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.04"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +4.18%  "

# Row 43
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +4.58%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.145.59"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.18%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "36.85"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.00%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.00"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.01%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.287.55"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.64%  "

# Row 48
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +3.37%  "

# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +7.01%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.68"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.64%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.08"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +3.01%  "
